# BIIBNamedTrade20.xlsx - "traded, fixed issues with the repeater"
#
# The trade-history "repeater" appends one new row per trade. This adds the
# next trade (row 8) below the existing rows (1 header + 6 trades, rows 1-7).
#
# Columns: A=Date, B=Profitable, C=Principle, D=Start Principle, E=BuyPrice,
#          F=SellPrice, G=IsShortSell, H=Price Change %, I=Strong trade

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the previous trade row's formatting down into the new row first (this
# carries over the date number-format on column A and the IsShortSell
# boolean style on column G, same as every other row in the repeater),
# then overwrite the copied values with the new trade's actual data.
$ws.Range("A7:I7").Copy($ws.Range("A8:I8"))

$ws.Cells.Item(8, 1).Value = 42654.746261574073
$ws.Cells.Item(8, 2).Value = $true
$ws.Cells.Item(8, 3).Value = 9836.08
$ws.Cells.Item(8, 4).Value = 9828.7099999999991
$ws.Cells.Item(8, 5).Value = 308
$ws.Cells.Item(8, 6).Value = 308.45999999999998
$ws.Cells.Item(8, 7).Value = $false
$ws.Cells.Item(8, 8).Value = 0.15
$ws.Cells.Item(8, 9).Value = $false
